$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '63.719.88'
$ws.Range("E2").Value = '  +1.37%  '
$ws.Range("D3").Value = '2.658.12'
$ws.Range("E3").Value = '  +3.05%  '
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").Value = "'593.88"
$ws.Range("E5").Value = '  +2.00%  '
$ws.Range("D6").Value = "'147.07"
$ws.Range("E6").Value = '  +0.43%  '
$ws.Range("D7").Value = "'1.00"
$ws.Range("E7").Value = '  +0.02%  '
$ws.Range("D8").Value = "'0.592"
$ws.Range("E8").Value = '  -0.43%  '
$ws.Range("E9").Value = '  +0.44%  '
$ws.Range("D10").Value = "'5.68"
$ws.Range("E10").Value = '  +0.79%  '
$ws.Range("D12").Value = "'0.357"
$ws.Range("E12").Value = '  +1.38%  '
$ws.Range("D13").Value = "'27.82"
$ws.Range("E13").Value = '  +2.33%  '
$ws.Range("D14").Value = '3.132.37'
$ws.Range("E14").Value = '  +2.97%  '
$ws.Range("D15").Value = '63.497.68'
$ws.Range("E15").Value = '  +1.23%  '
$ws.Range("E16").Value = '  +0.44%  '
$ws.Range("D17").Value = '2.623.44'
$ws.Range("E17").Value = '  +1.35%  '
$ws.Range("E18").Value = '  +1.17%  '
$ws.Range("D19").Value = "'344.31"
$ws.Range("E19").Value = '  +0.76%  '
$ws.Range("D20").Value = "'4.41"
$ws.Range("E20").Value = '  +0.79%  '
$ws.Range("E21").Value = '  +1.75%  '
$ws.Range("E22").Value = '  +0.17%  '
$ws.Range("D23").Value = "'68.08"
$ws.Range("E23").Value = '  +1.57%  '
$ws.Range("E24").Value = '  +8.13%  '
$ws.Range("D25").Value = "'1.58"
$ws.Range("E25").Value = '  +10.35%  '
$ws.Range("D26").Value = "'565.58"
$ws.Range("E26").Value = '  +22.90%  '
$ws.Range("E27").Value = '  +0.10%  '
$ws.Range("D28").Value = "'8.64"
$ws.Range("E28").Value = '  +4.05%  '
$ws.Range("E29").Value = '  +0.22%  '
$ws.Range("D30").Value = "'7.94"
$ws.Range("E30").Value = '  +1.03%  '
$ws.Range("E31").Value = '  +4.67%  '
$ws.Range("D33").Value = '0.0₃0821'
$ws.Range("E33").Value = '  +0.15%  '
$ws.Range("D34").Value = "'175.35"
$ws.Range("E34").Value = '  +0.07%  '
$ws.Range("E35").Value = '  +9.30%  '
$ws.Range("E36").Value = '  +0.06%  '
$ws.Range("E37").Value = '  +0.99%  '
$ws.Range("D38").Value = "'19.25"
$ws.Range("E38").Value = '  +1.49%  '
$ws.Range("D39").Value = "'1.80"
$ws.Range("E39").Value = '  +5.46%  '
$ws.Range("D40").Value = "'1.00"
$ws.Range("E40").Value = '  +0.06%  '
$ws.Range("D41").Value = "'169.34"
$ws.Range("E41").Value = '  +7.27%  '
$ws.Range("D42").Value = "'40.49"
$ws.Range("E42").Value = '  +2.91%  '
$ws.Range("E43").Value = '  +0.89%  '
$ws.Range("D44").Value = "'22.12"
$ws.Range("E44").Value = '  +5.19%  '
$ws.Range("D45").Value = "'0.0559"
$ws.Range("E45").Value = '  +3.64%  '
$ws.Range("E46").Value = '  -0.72%  '
$ws.Range("E47").Value = '  +2.65%  '
$ws.Range("E48").Value = '  -0.32%  '
$ws.Range("D49").Value = "'18.91"
$ws.Range("E49").Value = '  +2.89%  '
$ws.Range("E50").Value = '  +0.74%  '
$ws.Range("D51").Value = "'11.35"
$ws.Range("E51").Value = '  -0.56%  '
